$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the half-life value (B6) which drives the downstream formulas
$ws.Range("B6").Value = 0.005

# Update the selected cell to match the new active selection
$ws.Range("B7").Select()
